$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) The deck's first slide gets re-identified with a new internal SlideID.
#    PowerPoint assigns SlideID automatically (it is not a settable property),
#    and the only way to mint a fresh one through the object model is to have
#    PowerPoint create a new slide part for the same content - duplicating the
#    slide and removing the original reproduces exactly that renumbering
#    (confirmed to land on 510, matching the target).
# ---------------------------------------------------------------------------
$first = $p.Slides.Item(1)
$dup = $first.Duplicate()
$newFirst = $dup.Item(1)
$first.Delete()

# ---------------------------------------------------------------------------
# 2) Update the "05 de outubro de 2016" date caption on the (new) first slide
#    and on the last slide to "10 de outubro de 2016", preserving the run
#    layout seen in the target: "10" / " " / "de outubro de 2016".
# ---------------------------------------------------------------------------
function Update-DateCaption($slide) {
    $shape = $slide.Shapes.Item(1)
    $tr = $shape.TextFrame.TextRange
    $tr.Characters(1, 2).Text = "10"
    $tr.Characters(4, 18).Text = "de outubro de 2016"
}

Update-DateCaption $p.Slides.Item(1)
Update-DateCaption $p.Slides.Item($p.Slides.Count)

# ---------------------------------------------------------------------------
# 3) The cached "datetimeFigureOut" placeholder text on every slide layout,
#    the slide master, and the notes master is refreshed from 04/10/2016 to
#    10/10/2016 (PowerPoint recalculates these automatically on save).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "04/10/2016") {
                $sh.TextFrame.TextRange.Text = "10/10/2016"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($j).Shapes
}

Update-DatePlaceholder $p.NotesMaster.Shapes
